$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "292.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-6.60%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "40.34"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.01%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.011"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.27%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07321"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.34%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.523"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-8.85%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9265"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.18%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.378"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.90%"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.35%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1739"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.34%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04334"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.08%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08616"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.47%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1054"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.14%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001272"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.38%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005999"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.80%"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.33%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.295"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.19%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.972"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "4.99%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.01%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2795"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.50%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.03940"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.84%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.66%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003778"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-5.04%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001282"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.86%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003727"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-95.04%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02288"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-5.76%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.04982"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.22%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.005337"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "84.24%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007686"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.43%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1286"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.08%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007322"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.96%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007910"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.58%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3177"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.18%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006314"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-4.12%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.18%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.02045"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-92.40%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002103"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.18%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002003"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.18%"
